$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 180; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    $val = $cell.Value()
    if ($val -ne $null -and $val -like "/root/*") {
        $cell.Value = $val -replace "^/root/", "/dementia/"
    }
}

$ws.Names.Add("_xlnm._FilterDatabase_0", "=Sheet1!`$A`$1:`$M`$180")
